$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column (Price): force text representation via leading apostrophe, then
# reset the cell style back to Normal so no stray style index is introduced.
# E column (Volume): plain string assignment (Excel keeps these as text already).

$ws.Range("D2").Value = "'26.355.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.73%  '

$ws.Range("D3").Value = "'1.668.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.91%  '

$ws.Range("E4").Value = '  +0.53%  '

$ws.Range("D5").Value = "'220.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.48%  '

$ws.Range("D6").Value = "'0.5319"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.37%  '

$ws.Range("E7").Value = '  +0.48%  '

$ws.Range("D8").Value = "'0.2655"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.32%  '

$ws.Range("D9").Value = "'0.06362"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.60%  '

$ws.Range("D10").Value = "'20.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.05%  '

$ws.Range("D11").Value = "'0.07851"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.61%  '

$ws.Range("D12").Value = "'4.521"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.12%  '

$ws.Range("D13").Value = "'1.674.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.82%  '

$ws.Range("D14").Value = "'1.898.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.96%  '

$ws.Range("D15").Value = "'0.5593"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.87%  '

$ws.Range("D16").Value = "'0.0₅8171"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.14%  '

$ws.Range("D17").Value = "'66.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.01%  '

$ws.Range("D18").Value = "'26.382.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.90%  '

$ws.Range("E19").Value = '  +0.58%  '

$ws.Range("D20").Value = "'4.711"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.58%  '

$ws.Range("D21").Value = "'197.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.35%  '

$ws.Range("D22").Value = "'10.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.03%  '

$ws.Range("D23").Value = "'6.048"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.74%  '

$ws.Range("E24").Value = '  +0.45%  '

$ws.Range("D25").Value = "'145.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.15%  '

$ws.Range("D26").Value = "'0.1223"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("D27").Value = "'7.240"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.62%  '

$ws.Range("D28").Value = "'16.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.32%  '

$ws.Range("D29").Value = "'1.505"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.30%  '

$ws.Range("D30").Value = "'0.05906"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.97%  '

$ws.Range("D32").Value = "'3.557"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.22%  '

$ws.Range("D33").Value = "'3.321"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.77%  '

$ws.Range("D34").Value = "'1.605"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.92%  '

$ws.Range("D35").Value = "'0.9664"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.87%  '

$ws.Range("D36").Value = "'2.835"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.13%  '

$ws.Range("D37").Value = "'2.440"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.71%  '

$ws.Range("D38").Value = "'0.5815"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.60%  '

$ws.Range("D39").Value = "'0.01614"
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").Value = "'1.077.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.63%  '

$ws.Range("D41").Value = "'5.935"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").Value = "'0.8629"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.55%  '

$ws.Range("E43").Value = '  +0.55%  '

$ws.Range("D44").Value = "'102.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.15%  '

$ws.Range("D45").Value = "'1.810.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.96%  '

$ws.Range("D46").Value = "'58.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.00%  '

$ws.Range("D47").Value = "'0.0₈106"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.88%  '

$ws.Range("E48").Value = '  +0.74%  '

$ws.Range("D49").Value = "'0.4417"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.39%  '

$ws.Range("D50").Value = "'8.024"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.22%  '

$ws.Range("E51").Value = '  +0.08%  '

